# Modify the "Results of Testing" section (column G on the "Summary" sheet)
# so that it reflects the *area* in which confidence is increased, rather
# than referring to "confidence" directly.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Summary")

# --- Body rows (set first so new shared-string entries land in the same
#     order the original author's workbook ended up with) --------------
$ws.Range("G3").Value = "Module correctness"
$ws.Range("G4").Value = "Module interface correctness"
$ws.Range("G5").Value = "System correctness`nReliability"
$ws.Range("G6").Value = "System correctness`nUsability`nValidation"
$ws.Range("G7").Value = "Installability`nPortability"

# --- Header -------------------------------------------------------------
$ws.Range("G2").Value = "Results of Testing`n(Area of Confidence)"

# --- Cosmetic follow-up: Excel re-wraps/re-sizes row 2 & 3 and re-runs
#     "best fit" on column G once the text in those cells changes. The
#     row heights / column width below reproduce that recalculated
#     layout as closely as this host allows.
$ws.Rows.Item(2).RowHeight = 29.15
$ws.Rows.Item(3).RowHeight = 43.75
$ws.Columns.Item(7).ColumnWidth = 17.15

# --- Selection: the author had clicked into D4 by the time of the save
$ws.Activate()
$ws.Range("D4").Select()
